$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("D18").Value = 5

$ws.Range("D19").Select()
